$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Divide and conquer: replace the string node-label headers (row 1 and
# column A, formerly "EB", "EF I", "EF II", ... pulled from shared strings)
# with a numeric 1..8 index, generated by formulas that add 1 to the
# previous cell. ---

# Row 1: B1 is the seed (1), C1..I1 increment the cell to their left.
$ws.Range("B1").Value = 1
$ws.Range("C1").Formula = "=B1+1"
$ws.Range("D1:I1").FormulaR1C1 = "=RC[-1]+1"

# Column A: A2 is the seed (1), A3..A9 increment the cell above.
$ws.Range("A2").Value = 1
$ws.Range("A3").Formula = "=A2+1"
$ws.Range("A4:A9").FormulaR1C1 = "=R[-1]C+1"

# Move the active selection (cosmetic, matches the authored session state).
$ws.Range("J14").Select()
